$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 becomes the old Diodes (x4) / 1N4001 entry
$ws.Range("A10").Value = "Diodes (x4)"
$ws.Range("B10").Value = "1N4001"
$ws.Range("C10").Formula = "=0.21/10*4"
$ws.Range("D10").Value = "Futurlec"
$ws.Range("E10").Value = "Ok"

# Row 11 becomes the new MOSFETs (x4) / PMV31XN entry, received/ordered from RS Online
$ws.Range("A11").Value = "MOSFETs (x4)"
$ws.Range("B11").Value = "PMV31XN"
$ws.Range("C11").Formula = "=0.286*4"
$ws.Range("D11").Value = "RS Online"
$ws.Range("E11").Value = "Ordered"

# Motors (x4) received
$ws.Range("E12").Value = "Ok"

# Props received (came with motors)
$ws.Range("E17").Value = "Ok"

[void]$ws.Range("D11").Select()

$wb.Save()
